$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (General-looking numeric strings like "0.520" or "5.77" must
# stay literal text, matching the source inlineStr cells) by pre-setting
# the cell NumberFormat to "@" before assigning the Value.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.928.09'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.03%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.550.02'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.25%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.15%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '206.73'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.50%  '

$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.59%  '

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.12%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.05'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +3.02%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.246'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.21%  '

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.72%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0855'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.25%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.771.67'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.34%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.542.02'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.26%  '

$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.26%  '

$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.41%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '26.926.47'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.06%  '

$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.03%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '217.39'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.71%  '

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.71%  '

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.93%  '

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.18%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.04'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.62%  '

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.20%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.96'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.93%  '

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.47%  '

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.22%  '

$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.56%  '

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.97%  '

$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.13%  '

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.92%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.09'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.44%  '

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.07%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.420.24'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +3.96%  '

$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +3.97%  '

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +3.14%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.969'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.11%  '

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.17%  '

$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.43%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.520'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.59%  '

$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.77'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +5.35%  '

$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.807'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.22%  '

$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.16%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.31'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +4.20%  '

$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.48%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '64.28'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.35%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.684.98'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.33%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '87.65'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.59%  '

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.44%  '

$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +4.13%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0949'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.08%  '
